# Logs.xlsx visual/content correction
# - Re-sequences the Data/Description columns to reflect a consistent request
#   log (GET/POST/PUT/DELETE timestamps), matching the new log capture.
# - Makes the Method-colour fonts bold and recolours them to be readable
#   (cyan -> black, spring-green -> dark-green, violet -> navy) and whitens
#   the row fill (was mid-grey).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Cell content (columns B/C/D, rows 2-15)
# ---------------------------------------------------------------------
$rows = @(
    @{ Row=2;  B='GET';    C="user 'admin' gets all products from database";        D='2022-12-14 09:25:06' },
    @{ Row=3;  B='POST';   C="user 'admin' inserted product '4' from database";     D='2022-12-14 09:25:13' },
    @{ Row=4;  B='GET';    C="user 'admin' gets product '4' from database";         D='2022-12-14 09:25:13' },
    @{ Row=5;  B='GET';    C="user 'admin' gets all products from database";        D='2022-12-14 09:25:13' },
    @{ Row=6;  B='GET';    C="user 'admin' gets product '4' from database";         D='2022-12-14 09:25:18' },
    @{ Row=7;  B='PUT';    C="user 'admin' updated product '4' from database";      D='2022-12-14 09:25:21' },
    @{ Row=8;  B='GET';    C="user 'admin' gets product '4' from database";         D='2022-12-14 09:25:21' },
    @{ Row=9;  B='GET';    C="user 'admin' gets all products from database";        D='2022-12-14 09:25:21' },
    @{ Row=10; B='GET';    C="user 'admin' gets product '4' from database";         D='2022-12-14 09:25:22' },
    @{ Row=11; B='PUT';    C="user 'admin' updated product '4' from database";      D='2022-12-14 09:25:25' },
    @{ Row=12; B='GET';    C="user 'admin' gets product '4' from database";         D='2022-12-14 09:25:25' },
    @{ Row=13; B='GET';    C="user 'admin' gets all products from database";        D='2022-12-14 09:25:25' },
    @{ Row=14; B='DELETE'; C="user 'admin' deleted product '4' from database";      D='2022-12-14 09:25:27' },
    @{ Row=15; B='GET';    C="user 'admin' gets all products from database";        D='2022-12-14 09:25:27' }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
}

# ---------------------------------------------------------------------
# 2. Row highlight colour, now tied to the (possibly moved) Method value
# ---------------------------------------------------------------------
# style groups, expressed as contiguous row ranges sharing a Method
$putRows    = @(7, 11)
$postRows   = @(3)
$deleteRows = @(14)
# every other data row (2,4,5,6,8,9,10,12,13,15) keeps the GET look

foreach ($n in $postRows) {
    $rng = $ws.Range("A${n}:D${n}")
    $rng.Font.Bold = $true
    $rng.Font.Color = 25600       # dark green FF006400
    $rng.Interior.Color = 16777215 # white FFFFFFFF
}

foreach ($n in $putRows) {
    $rng = $ws.Range("A${n}:D${n}")
    $rng.Font.Bold = $true
    $rng.Font.Color = 10040064    # navy blue FF003399
    $rng.Interior.Color = 16777215 # white FFFFFFFF
}

foreach ($n in $deleteRows) {
    $rng = $ws.Range("A${n}:D${n}")
    $rng.Font.Bold = $true
    $rng.Font.Color = 2237106     # firebrick FFB22222 (unchanged colour, now bold)
    $rng.Interior.Color = 16777215 # white FFFFFFFF
}

# GET rows (everything not POST/PUT/DELETE): bold black text, white fill
$getRowRanges = @("A2:D2", "A4:D6", "A8:D10", "A12:D13", "A15:D15")
foreach ($addr in $getRowRanges) {
    $rng = $ws.Range($addr)
    $rng.Font.Bold = $true
    $rng.Font.Color = 0            # black FF000000
    $rng.Interior.Color = 16777215 # white FFFFFFFF
}
